$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, shifting existing rows 17-92 down to 18-93
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with fresh data
$ws.Cells.Item(17, 1).Value = 8
$ws.Cells.Item(17, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(17, 3).Value = "Coquimbo"
$ws.Cells.Item(17, 4).Value = 44453
$ws.Cells.Item(17, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17, 5).Value = 4
$ws.Cells.Item(17, 6).Value = 100112021
$ws.Cells.Item(17, 7).Value = "Ají"
$ws.Cells.Item(17, 8).Value = "Inferno"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 600
$ws.Cells.Item(17, 11).Value = 38000
$ws.Cells.Item(17, 12).Value = 39000
$ws.Cells.Item(17, 13).Value = 38500
$ws.Cells.Item(17, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(17, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(17, 16).Value = 3208
$ws.Cells.Item(17, 17).Value = 12
$ws.Cells.Item(17, 18).Value = "Hortaliza"
